$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resize the table / list object to cover the new data range A1:B7
$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:B7"))

# Rewrite the data rows with the new, more realistic demo values
$ws.Range("A2").Value = "personal_data"
$ws.Range("B2").Value = "bevnat"

$ws.Range("A3").Value = "personal_data"
$ws.Range("B3").Value = "statpop"

$ws.Range("A4").Value = "sensible_data"
$ws.Range("B4").Value = "bevnat"

$ws.Range("A5").Value = "population"
$ws.Range("B5").Value = "bevnat"

$ws.Range("A6").Value = "population"
$ws.Range("B6").Value = "statpop"

$ws.Range("A7").Value = "societe"
$ws.Range("B7").Value = "statpop"

# Widen column A to fit the longer labels now used
$ws.Columns.Item(1).ColumnWidth = 15.83

# Move the active selection to reflect where the author left off editing
$ws.Range("B9").Select() | Out-Null
